# Remove the blank paragraph, the "Ver no Jupiter..." paragraph, the
# following blank paragraph, and the trailing page-break paragraph that
# used to sit right after the "LOT2013: ..." requisito paragraph.
$d = $word.ActiveDocument

function Find-ParagraphIndexByText($searchText) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $found = $r.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                              $true, 1, $false, "", 0)
    if (-not $found) {
        return -1
    }
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {
            return $i
        }
    }
    return -1
}

$lot2013Index = Find-ParagraphIndexByText("LOT2013: Engenharia Bioqu" + [char]237 + "mica I (Requisito fraco)")
$jupiterIndex = Find-ParagraphIndexByText("Ver no Jupiter Salvar em pdf Salvar em docx")

# The four paragraphs to delete: the blank line right after LOT2013, the
# "Ver no Jupiter..." line itself, the blank line right after it, and the
# following page-break paragraph.
$startPara = $d.Paragraphs.Item($lot2013Index + 1)
$endPara = $d.Paragraphs.Item($jupiterIndex + 2)

$range = $d.Range($startPara.Range.Start, $endPara.Range.End)
$range.Delete()
